$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2892.1667
$ws.Range("I113").Value = 2213.25
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 2213.25
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = 1040.75
$ws.Range("N113").Value = -10758
$ws.Range("H137").Value = 2794.5386
$ws.Range("I137").Value = 1921.5834
$ws.Range("J137").Value = 3542.7856
$ws.Range("K137").Value = 5764.7502
$ws.Range("L137").Value = 10628.3568
$ws.Range("M137").Value = -3214.7502
$ws.Range("N137").Value = -15728.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2918.1965
$ws.Range("I32").Value = 3132.319
$ws.Range("K32").Value = 3132.319
$ws.Range("M32").Value = -2845.319
$ws.Range("H80").Value = 4001
$ws.Range("I80").Value = 4001
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4001
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3003
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 4001
$ws.Range("I83").Value = 4001
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12003
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7011
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 600
$ws.Range("I25").Value = 600
$ws.Range("K25").Value = 600
$ws.Range("M25").Value = -365
$ws.Range("H82").Value = 24042.5
$ws.Range("I82").Value = 21651
$ws.Range("K82").Value = 21651
$ws.Range("M82").Value = -21268
$ws.Range("H85").Value = 24042.5
$ws.Range("I85").Value = 21651
$ws.Range("K85").Value = 21651
$ws.Range("M85").Value = -20325
$ws.Range("H134").Value = 3531.4102
$ws.Range("I134").Value = 855.5
$ws.Range("K134").Value = 2566.5
$ws.Range("M134").Value = -31.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1740.5745
$ws.Range("J31").Value = 2055.3794
$ws.Range("L31").Value = 2055.3794
$ws.Range("N31").Value = -2645.3794
$ws.Range("H34").Value = 1740.5745
$ws.Range("J34").Value = 2055.3794
$ws.Range("L34").Value = 2055.3794
$ws.Range("N34").Value = -2459.3794
$ws.Range("H58").Value = 700.7917
$ws.Range("I58").Value = 664.3684
$ws.Range("J58").Value = 839.2
$ws.Range("K58").Value = 664.3684
$ws.Range("L58").Value = 839.2
$ws.Range("M58").Value = -461.3684
$ws.Range("N58").Value = -1245.2
$ws.Range("H74").Value = 30078.25
$ws.Range("J74").Value = 30078.25
$ws.Range("L74").Value = 30078.25
$ws.Range("N74").Value = -31826.25
$ws.Range("H77").Value = 30078.25
$ws.Range("J77").Value = 30078.25
$ws.Range("L77").Value = 90234.75
$ws.Range("N77").Value = -98970.75
$ws.Range("H86").Value = 4779850
$ws.Range("I86").Value = 6670347.5
$ws.Range("K86").Value = 6670347.5
$ws.Range("M86").Value = -6669224.5
$ws.Range("H89").Value = 4779850
$ws.Range("I89").Value = 6670347.5
$ws.Range("K89").Value = 33351737.5
$ws.Range("M89").Value = -33346121.5
$ws.Range("H136").Value = 700.7917
$ws.Range("I136").Value = 664.3684
$ws.Range("J136").Value = 839.2
$ws.Range("K136").Value = 1993.1052
$ws.Range("L136").Value = 2517.6
$ws.Range("M136").Value = 556.8948
$ws.Range("N136").Value = -7617.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1666.6666
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -4331
$ws.Range("N17").Value = -6338
$ws.Range("H39").Value = 1407.9231
$ws.Range("J39").Value = 1491.6666
$ws.Range("L39").Value = 4474.9998
$ws.Range("N39").Value = -5062.9998
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 2597.4534
$ws.Range("I68").Value = 1101
$ws.Range("J68").Value = 2709.6875
$ws.Range("K68").Value = 3303
$ws.Range("L68").Value = 8129.0625
$ws.Range("M68").Value = -2492
$ws.Range("N68").Value = -9751.0625
$ws.Range("H71").Value = 2597.4534
$ws.Range("I71").Value = 1101
$ws.Range("J71").Value = 2709.6875
$ws.Range("K71").Value = 9909
$ws.Range("L71").Value = 24387.1875
$ws.Range("M71").Value = -5853
$ws.Range("N71").Value = -32499.1875
$ws.Range("H86").Value = 1352
$ws.Range("I86").Value = 1352
$ws.Range("K86").Value = 4056
$ws.Range("M86").Value = -2870
$ws.Range("H89").Value = 1352
$ws.Range("I89").Value = 1352
$ws.Range("K89").Value = 12168
$ws.Range("M89").Value = -6240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4896.8
$ws.Range("I80").Value = 2840
$ws.Range("J80").Value = 6953.6
$ws.Range("K80").Value = 2840
$ws.Range("L80").Value = 6953.6
$ws.Range("M80").Value = -1842
$ws.Range("N80").Value = -8949.6
$ws.Range("H83").Value = 4896.8
$ws.Range("I83").Value = 2840
$ws.Range("J83").Value = 6953.6
$ws.Range("K83").Value = 14200
$ws.Range("L83").Value = 34768
$ws.Range("M83").Value = -9208
$ws.Range("N83").Value = -44752
$ws.Range("H97").Value = 2331.3845
$ws.Range("I97").Value = 2400.6667
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 2400.6667
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -1904.6667
$ws.Range("N97").Value = -2492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 987.57574
$ws.Range("I61").Value = 945.38464
$ws.Range("J61").Value = 1144.2858
$ws.Range("K61").Value = 945.38464
$ws.Range("L61").Value = 1144.2858
$ws.Range("M61").Value = -743.38464
$ws.Range("N61").Value = -1548.2858
$ws.Range("H113").Value = 987.57574
$ws.Range("I113").Value = 945.38464
$ws.Range("J113").Value = 1144.2858
$ws.Range("K113").Value = 945.38464
$ws.Range("L113").Value = 1144.2858
$ws.Range("M113").Value = 1224.61536
$ws.Range("N113").Value = -5484.2858
$ws.Range("H136").Value = 17900.5
$ws.Range("I136").Value = 20880.6
$ws.Range("K136").Value = 62641.8
$ws.Range("M136").Value = -60091.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7811.5
$ws.Range("I81").Value = 622.5
$ws.Range("J81").Value = 15000.5
$ws.Range("K81").Value = 1245
$ws.Range("L81").Value = 30001
$ws.Range("M81").Value = -184
$ws.Range("N81").Value = -32123
$ws.Range("H84").Value = 7811.5
$ws.Range("I84").Value = 622.5
$ws.Range("J84").Value = 15000.5
$ws.Range("K84").Value = 6225
$ws.Range("L84").Value = 150005
$ws.Range("M84").Value = -921
$ws.Range("N84").Value = -160613
$ws.Range("H86").Value = 7325
$ws.Range("J86").Value = 7325
$ws.Range("L86").Value = 7325
$ws.Range("N86").Value = -9571
$ws.Range("H89").Value = 7325
$ws.Range("J89").Value = 7325
$ws.Range("L89").Value = 36625
$ws.Range("N89").Value = -47857